$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 204
$ws.Range("D21").Value = 182
$ws.Range("E21").Value = 22
$ws.Range("F21").Value = 52.14899713467048
